$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 4)
$ws.Range("A4").Value = "Two Pointers"
$ws.Range("B4").Value = "Valid Palindrome"
$ws.Range("C4").Value = "two pointers, string;"

# Copy the style from B3 (e.g. special font color) onto B4 to match existing formatting
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to the new last cell, matching the saved view state
$ws.Range("C4").Select()
